# QA Excel Compiler - STATUS tracking update
# Normalizes old STATUS values (OK/ERROR) to the new valid set
# (ISSUE, NO ISSUE, BLOCKED) and trims a couple of COMMENT strings,
# per the mock test file updates described in the commit message.

$wb = $excel.ActiveWorkbook

# --- Sheet1 ---
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("E2").Value = "NO ISSUE"

$ws1.Range("E3").Value = "NO ISSUE"
$ws1.Range("F3").Value = "I think it's fine"

$ws1.Range("E4").Value = "ISSUE"
$ws1.Range("F4").Value = "Missing article"

# --- Sheet2 ---
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("E3").Value = "NO ISSUE"
$ws2.Range("F3").Value = "Shop verified"
